# This script reproduces the commit:
#   "cambio unidades reducción, implementación R2, renaming módulos"
#
# Summary of the change (derived from the OOXML diff):
#  1) The 5 worksheets are renamed from V2_0.002 / V2_0.0015 / V2_0.001 /
#     V2_0.0005 / V2_0  to  V2_0 / V2_0.25 / V2_0.5 / V2_0.75 / V2_1
#     (same relationship ids / same physical sheets, just a rename so the
#     tab names represent the fraction of reduction applied: 0, .25, .5, .75, 1).
#  2) For every data row (2..88) the first sheet (fraction = 0) holds the
#     baseline "Huella / tCO2e" value in column C and baseline
#     "Huella/cápita / tCO2e" in column E. Column D held the reduction
#     fraction itself (0, 0.25, 0.5, 0.75, 1) before the edit.
#     After the edit:
#       - column C (Huella / tCO2e)          = baselineC * (1-frac) * 0.2
#       - column D (Reducción / %  -> now "reduced" footprint kept)
#                                              = baselineC * (0.8 + 0.2*frac)
#       - column E (Huella/cápita / tCO2e)    = baselineE * (1-frac) * 0.2
#     i.e. values now reflect a R2 recalculation where only 20% of the
#     original footprint is considered reducible and the reduction is
#     redistributed between "remaining" (C) and "reduced" (D) columns.

$wb = $excel.ActiveWorkbook

# --- worksheet handles (order on tabs before the rename) -------------------
$ws1 = $wb.Worksheets.Item(1)   # was "V2_0.002"   -> becomes "V2_0"     (frac 0)
$ws2 = $wb.Worksheets.Item(2)   # was "V2_0.0015"  -> becomes "V2_0.25"  (frac .25)
$ws3 = $wb.Worksheets.Item(3)   # was "V2_0.001"   -> becomes "V2_0.5"   (frac .5)
$ws4 = $wb.Worksheets.Item(4)   # was "V2_0.0005"  -> becomes "V2_0.75"  (frac .75)
$ws5 = $wb.Worksheets.Item(5)   # was "V2_0"       -> becomes "V2_1"     (frac 1)

$firstRow = 2
$lastRow = 88

# --- capture the original (baseline) column C & E values from sheet 1 ------
# (must be read before any sheet is modified, sheet1 holds the source data)
$baseC = @{}
$baseE = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $baseC[$r] = $ws1.Cells.Item($r, 3).Value()
    $baseE[$r] = $ws1.Cells.Item($r, 5).Value()
}

# --- apply the new values to every sheet, using its reduction fraction -----
$sheetsAndFracs = @(
    @{ Sheet = $ws1; Frac = 0.0 },
    @{ Sheet = $ws2; Frac = 0.25 },
    @{ Sheet = $ws3; Frac = 0.5 },
    @{ Sheet = $ws4; Frac = 0.75 },
    @{ Sheet = $ws5; Frac = 1.0 }
)

foreach ($entry in $sheetsAndFracs) {
    $ws = $entry.Sheet
    $frac = $entry.Frac
    for ($r = $firstRow; $r -le $lastRow; $r++) {
        $c1 = $baseC[$r]
        $e1 = $baseE[$r]

        $newC = $c1 * (1 - $frac) * 0.2
        $newD = $c1 * (0.8 + 0.2 * $frac)
        $newE = $e1 * (1 - $frac) * 0.2

        $ws.Cells.Item($r, 3).Value = $newC
        $ws.Cells.Item($r, 4).Value = $newD
        $ws.Cells.Item($r, 5).Value = $newE
    }
}

# --- rename the worksheet tabs ---------------------------------------------
# Rename sheet 5 first to avoid a transient name collision with sheet 1's
# target name ("V2_0", which is sheet 5's current name).
$ws5.Name = "V2_1"
$ws1.Name = "V2_0"
$ws2.Name = "V2_0.25"
$ws3.Name = "V2_0.5"
$ws4.Name = "V2_0.75"
